$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Cronograma #1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cronograma #1")

# Text updates on existing rows (above the insertion point, row numbers unaffected)
$ws.Range("B54").Value = "Generar documentación para el usuario versión Preliminar"
$ws.Range("D54").Value = "SGI-MU1.DOCX"
$ws.Range("E57").Value = "Balarezo/JP, Canecillas/A, Balceda/PB"

# Mark the now-completed tasks as 100% done
$ws.Range("H53").Value = 1
$ws.Range("H54").Value = 1
$ws.Range("H55").Value = 1
$ws.Range("H56").Value = 1
$ws.Range("H57").Value = 1
$ws.Range("H58").Value = 1

# Insert a new row for the "versión Final" user-documentation deliverable,
# right before the "Documentar Acta de finalización del Hito 3" row.
$ws.Rows("59").Insert(-4121, 0)
$ws.Rows("59").RowHeight = 15.75

$ws.Range("B59").Value = "Generar documentación para el usuario versión Versión Final"
$ws.Range("C59").Value = "Manual de usuario"
$ws.Range("D59").Value = "SGI-MU2.DOCX"
$ws.Range("E59").Value = "Soller/PB, Balceda/PB, Huarhua/UI, Del Aguila/QA"
$ws.Range("F59").Value = 45103
$ws.Range("G59").Value = 45107
$ws.Range("H59").Value = 1

# Rows pushed down by the insert keep their old content but need refreshed
# completion percentages and one author swap.
$ws.Range("H60").Value = 1
$ws.Range("E60").Value = "Justiniano/A, Del Aguila/QA, Soller/PB"
$ws.Range("H61").Value = 1
$ws.Range("H62").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Control de Versiones"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Control de Versiones")

$ws3.Range("E7").Value = ""

$ws3.Range("A8").Value = 44960
$ws3.Range("B8").Value = 45107
$ws3.Range("C8").Value = "Luis Balarezo"
$ws3.Range("D8").Value = "Actualización del avance de los ítems al 100%"
